$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings in column D stay as text, matching source formatting
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '61.256.25'
$ws.Range("E2").Value = '  -2.12%  '
$ws.Range("D3").Value = '3.015.88'
$ws.Range("E3").Value = '  -4.53%  '
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").Value = '569.10'
$ws.Range("E5").Value = '  -3.35%  '
$ws.Range("D6").Value = '129.16'
$ws.Range("E6").Value = '  -4.97%  '
$ws.Range("E7").Value = '  -0.17%  '
$ws.Range("D8").Value = '3.016.20'
$ws.Range("E8").Value = '  -4.25%  '
$ws.Range("D9").Value = '0.497'
$ws.Range("E9").Value = '  -2.01%  '
$ws.Range("D10").Value = '0.135'
$ws.Range("E10").Value = '  -4.91%  '
$ws.Range("D11").Value = '5.18'
$ws.Range("E11").Value = '  -1.41%  '
$ws.Range("E12").Value = '  -5.14%  '
$ws.Range("D13").Value = '0.0000224'
$ws.Range("E13").Value = '  -4.46%  '
$ws.Range("D14").Value = '32.90'
$ws.Range("E14").Value = '  -2.43%  '
$ws.Range("E15").Value = '  -0.47%  '
$ws.Range("D16").Value = '3.511.10'
$ws.Range("E16").Value = '  -4.64%  '
$ws.Range("D17").Value = '61.200.61'
$ws.Range("E17").Value = '  -2.25%  '
$ws.Range("D18").Value = '3.009.85'
$ws.Range("E18").Value = '  -4.86%  '
$ws.Range("D19").Value = '6.23'
$ws.Range("E19").Value = '  -4.87%  '
$ws.Range("D20").Value = '439.51'
$ws.Range("E20").Value = '  -3.08%  '
$ws.Range("D21").Value = '13.19'
$ws.Range("E21").Value = '  -5.60%  '
$ws.Range("D22").Value = '0.662'
$ws.Range("E22").Value = '  -5.48%  '
$ws.Range("D23").Value = '7.15'
$ws.Range("E23").Value = '  -5.96%  '
$ws.Range("B24").Value = 'Litecoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D24").Value = '78.98'
$ws.Range("E24").Value = '  -5.49%  '
$ws.Range("B25").Value = 'InternetComputer(DFINITY)'
$ws.Range("C25").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D25").Value = '12.55'
$ws.Range("E25").Value = '  -5.95%  '
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  +0.09%  '
$ws.Range("D27").Value = '0.999'
$ws.Range("E27").Value = '  -0.18%  '
$ws.Range("D28").Value = '2.50'
$ws.Range("E28").Value = '  -6.87%  '
$ws.Range("D29").Value = '7.24'
$ws.Range("E29").Value = '  -6.14%  '
$ws.Range("D30").Value = '6.22'
$ws.Range("E30").Value = '  -7.72%  '
$ws.Range("D31").Value = '1.88'
$ws.Range("E31").Value = '  -6.82%  '
$ws.Range("D32").Value = '25.54'
$ws.Range("E32").Value = '  -6.16%  '
$ws.Range("D33").Value = '0.0944'
$ws.Range("E33").Value = '  -8.28%  '
$ws.Range("E34").Value = '  -4.83%  '
$ws.Range("D35").Value = '0.956'
$ws.Range("E35").Value = '  -7.19%  '
$ws.Range("D36").Value = '5.57'
$ws.Range("E36").Value = '  -4.21%  '
$ws.Range("D37").Value = '50.10'
$ws.Range("E37").Value = '  -2.09%  '
$ws.Range("D38").Value = '0.0₃0683'
$ws.Range("E38").Value = '  -3.10%  '
$ws.Range("D39").Value = '0.0363'
$ws.Range("E39").Value = '  -6.56%  '
$ws.Range("D40").Value = '7.74'
$ws.Range("E40").Value = '  -3.60%  '
$ws.Range("E41").Value = '  -2.72%  '
$ws.Range("D42").Value = '373.82'
$ws.Range("E42").Value = '  -7.10%  '
$ws.Range("D43").Value = '2.655.09'
$ws.Range("E43").Value = '  -4.67%  '
$ws.Range("D44").Value = '2.44'
$ws.Range("E44").Value = '  -8.95%  '
$ws.Range("D46").Value = '0.236'
$ws.Range("E46").Value = '  -5.16%  '
$ws.Range("D47").Value = '119.53'
$ws.Range("E47").Value = '  -4.60%  '
$ws.Range("D48").Value = '1.97'
$ws.Range("E48").Value = '  -7.15%  '
$ws.Range("D49").Value = '33.08'
$ws.Range("E49").Value = '  -4.92%  '
$ws.Range("E50").Value = '  -4.21%  '
$ws.Range("D51").Value = '23.65'
$ws.Range("E51").Value = '  -6.55%  '
